$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 3429
$ws.Range("F6").Value = 4936
$ws.Range("F7").Value = 492
$ws.Range("F8").Value = 323
$ws.Range("F10").Value = 649
$ws.Range("F11").Value = 290
$ws.Range("F13").Value = 23
$ws.Range("F14").Value = 677
$ws.Range("F15").Value = 298
$ws.Range("F18").Value = 152
$ws.Range("F20").Value = 352
$ws.Range("F21").Value = 4808
$ws.Range("F25").Value = 5943
$ws.Range("F26").Value = 19
$ws.Range("F27").Value = 8
$ws.Range("F28").Value = 3205
$ws.Range("F29").Value = 289
$ws.Range("F30").Value = 688
$ws.Range("F31").Value = 4429
$ws.Range("F32").Value = 311
$ws.Range("F33").Value = 104
$ws.Range("F34").Value = 135
$ws.Range("F35").Value = 913
$ws.Range("F37").Value = 17
$ws.Range("F39").Value = 825
$ws.Range("F40").Value = 905

# Sheet 2: 演出 (Performance)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 16
$ws.Range("F6").Value = 52

# Sheet 3: 本地生活 (Local life)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 44

# Sheet 4: 全部类型 (All types)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F5").Value = 44
$ws.Range("F8").Value = 3429
$ws.Range("F10").Value = 4936
$ws.Range("F11").Value = 492
$ws.Range("F12").Value = 323
$ws.Range("F14").Value = 649
$ws.Range("F15").Value = 290
$ws.Range("F17").Value = 23
$ws.Range("F18").Value = 677
$ws.Range("F19").Value = 298
$ws.Range("F23").Value = 152
$ws.Range("F25").Value = 352
$ws.Range("F26").Value = 4808
$ws.Range("F30").Value = 5943
$ws.Range("F31").Value = 19
$ws.Range("F32").Value = 8
$ws.Range("F33").Value = 3205
$ws.Range("F34").Value = 289
$ws.Range("F35").Value = 688
$ws.Range("F36").Value = 4429
$ws.Range("F37").Value = 311
$ws.Range("F38").Value = 16
$ws.Range("F39").Value = 104
$ws.Range("F40").Value = 913
$ws.Range("F42").Value = 17
$ws.Range("F44").Value = 825
$ws.Range("F45").Value = 905
$ws.Range("F49").Value = 52
